# #5: property boat&car done
# The "汽車" (car) sheet (sheet3) previously had its header row (row 1)
# filled in with duplicate data values instead of real column headers.
# This fixes row 1 to use the same header schema as the other property
# sheets, and extends row 2 with the common trailer columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that were missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Row 1: proper column headers (style already applied from existing cells)
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Row 2: fill in the trailer columns that were missing from the data row.
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2011-11-21"
$ws.Range("K2").Value = "潘維剛"
$ws.Range("L2").Value = 678
$ws.Range("M2").Value = "tmpcafb1"
$ws.Range("N2").Value = 29
